$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 42647.681840277779
$ws.Range("B5").Value = $true
$ws.Range("C5").Value = 10104.16
$ws.Range("D5").Value = 10020.99
$ws.Range("E5").Value = 18.12
$ws.Range("F5").Value = 17.97
$ws.Range("G5").Value = $true
$ws.Range("H5").Value = -0.83
$ws.Range("I5").Value = $true

$ws.Range("A5").NumberFormat = "m/d/yy h:mm"
$ws.Range("G5").NumberFormat = "m/d/yy h:mm"
